# Update gh-pages output (Beijing comic-convention info) to the scrape
# generated at commit 456a3b4.
#
# Sheet "展览" (Exhibitions): rows 14-19 cascade - one finished event
# (MQ&THEBONE) drops off, a new event (不舍昼夜3.0-奇妙童话夜) appears,
# and the rows in between shift up by one with a couple of refreshed
# "want-to-go" counts (column F). Several other rows across the sheet
# just get their column-F attendance counter bumped.
# Sheet "演出" (Performances): one attendance counter bump.
# Sheet "全部类型" (All types, a superset of every category) already
# contains every event row, so it only needs the same column-F counter
# bumps - no row content changes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F9").Value  = 286
$ws.Range("F10").Value = 410
$ws.Range("F12").Value = 1812
$ws.Range("F13").Value = 792

# Row 14: was 北京·MQ&THEBONE首届怀旧only -> becomes 北京·万达国潮动漫嘉年华
$ws.Range("C14").Value = "北京·万达国潮动漫嘉年华【免票活动】"
$ws.Range("D14").Value = "新华西街58号万达广场 北京通州万达广场"
$ws.Range("E14").Value = "2024.09.15 13:00-09.15 18:00"
$ws.Range("F14").Value = 11
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=91479"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202408/84Z3tWJF1724833337880.jpeg"

# Row 15: was 北京·万达国潮动漫嘉年华 -> becomes 北京·乐多港万达中秋动漫次元嘉年华
$ws.Range("C15").Value = "北京·乐多港万达中秋动漫次元嘉年华【免票活动】"
$ws.Range("D15").Value = "城南街道南口路29号 北京乐多港万达广场"
$ws.Range("E15").Value = "2024.09.15 13:30-09.15 18:30"
$ws.Range("F15").Value = 14
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=91481"
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202408/UrBQ6ywo1724817393278.jpeg"

# Row 16: was 北京·乐多港万达中秋动漫次元嘉年华 -> becomes 北京·原神only4.0同人展
$ws.Range("C16").Value = "北京·原神only4.0同人展"
$ws.Range("D16").Value = "北花园路1号 超级蜂巢"
$ws.Range("E16").Value = "2024.09.15 10:00-09.15 17:00"
$ws.Range("F16").Value = 1553
$ws.Range("G16").Value = 68
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=87564"
$ws.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202407/EfEAeJDS1720776874376.jpeg"

# Row 17: stays 北京·原神only4.0同人展, only the counter bumps
$ws.Range("F17").Value = 1554

# Row 18: was 北京·原神only4.0同人展 -> becomes 北京·AINI二次元派对
# Force text so "2024-09-16" isn't auto-converted into a date serial by
# Excel's type inference (matches how the rest of column B is stored).
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "2024-09-16"
$ws.Range("C18").Value = "北京·AINI二次元派对【免票展会】"
$ws.Range("D18").Value = "新村街道丰科路6号F1-102-103 万达广场(丰科店)"
$ws.Range("E18").Value = "2024.09.16 10:00-09.16 18:00"
$ws.Range("F18").Value = 1292
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=90730"
$ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202408/9SUINRO61723558972754.jpeg"

# Row 19: was 北京·AINI二次元派对 -> becomes the new event 北京·不舍昼夜3.0-奇妙童话夜
$ws.Range("C19").Value = "北京·不舍昼夜3.0-奇妙童话夜"
$ws.Range("D19").Value = "酒仙桥北路2号院798艺术区706后街1号 北京格瑞斯艺术酒店"
$ws.Range("E19").Value = "2024.09.16 20:00-09.17 02:00"
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 158
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=91042"
$ws.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202408/ZRwTjxgi1724204402969.jpeg"

$ws.Range("F22").Value = 375
$ws.Range("F26").Value = 6792
$ws.Range("F27").Value = 7220
$ws.Range("F28").Value = 14
$ws.Range("F29").Value = 159
$ws.Range("F33").Value = 1
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = 32
$ws.Range("F37").Value = 1330
$ws.Range("F40").Value = 639
$ws.Range("F42").Value = 1341
$ws.Range("F43").Value = 280
$ws.Range("F44").Value = 155
$ws.Range("F45").Value = 163

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 19

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F12").Value = 286
$ws4.Range("F14").Value = 410
$ws4.Range("F16").Value = 1812
$ws4.Range("F17").Value = 792
$ws4.Range("F18").Value = 11
$ws4.Range("F19").Value = 14
$ws4.Range("F20").Value = 1554
$ws4.Range("F21").Value = 1554
$ws4.Range("F26").Value = 375
$ws4.Range("F30").Value = 6792
$ws4.Range("F31").Value = 7220
$ws4.Range("F33").Value = 1330
$ws4.Range("F39").Value = 639
$ws4.Range("F43").Value = 1341
$ws4.Range("F44").Value = 280
